$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe out the existing H:S block (values, styles, hyperlinks) ---
# so it can be rebuilt cleanly in the new column order.
$ws.Range("H1:S3").Hyperlinks.Delete()
$ws.Range("H1:S3").Clear()

# --- Row 1: headers (plain text, no hyperlink styling) ---
$ws.Range("H1").Value = "年"
$ws.Range("I1").Value = "アイテムURL"
$ws.Range("J1").Value = "IIIFマニフェストURI"
$ws.Range("K1").Value = "viewingDirection"
$ws.Range("L1").Value = "帰属"
$ws.Range("M1").Value = "ID"
$ws.Range("N1").Value = "ソート用項目"
$ws.Range("O1").Value = "コレクション"
$ws.Range("P1").Value = "機械可読ドキュメント"
$ws.Range("Q1").Value = "サムネイル"
$ws.Range("R1").Value = "ウェブサイトURL"
$ws.Range("S1").Value = "利用条件"

function Add-HyperlinkCell($cellRef, $address, $subAddress, $displayText) {
    $rng = $ws.Range($cellRef)
    if ($subAddress) {
        $ws.Hyperlinks.Add($rng, $address, $subAddress)
    } else {
        $ws.Hyperlinks.Add($rng, $address)
    }
    $rng.Value = $displayText
    # Match the workbook's existing hyperlink style (blue/underline) instead
    # of letting Excel auto-create a brand-new "Hyperlink" cell style.
    $rng.Font.Underline = $true
    $rng.Font.Color = 16711680
}

# --- Row 2 (item: ecd1285a-...) ---
$ws.Range("L2").Value = "東京大学総合図書館 / General Library in the University of Tokyo"
$ws.Range("M2").Value = "ecd1285a-42b6-4541-b640-1067b2f9fde0"
$ws.Range("O2").Value = "キリシタン写本"

Add-HyperlinkCell "I2" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian/document/ecd1285a-42b6-4541-b640-1067b2f9fde0" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian/document/ecd1285a-42b6-4541-b640-1067b2f9fde0"
Add-HyperlinkCell "J2" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/39287/manifest" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/39287/manifest"
Add-HyperlinkCell "K2" "http://iiif.io/api/presentation/2" "rightToLeftDirection" "http://iiif.io/api/presentation/2#rightToLeftDirection"
Add-HyperlinkCell "P2" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/39287" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/39287"
Add-HyperlinkCell "Q2" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/files/square/ebb775777d3b47a521407bbed523d5f088ba43f0.jpg" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/files/square/ebb775777d3b47a521407bbed523d5f088ba43f0.jpg"
Add-HyperlinkCell "R2" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian"
Add-HyperlinkCell "S2" "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse" $null "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse"

# --- Row 3 (item: 394c6988-...) ---
$ws.Range("L3").Value = "東京大学総合図書館 / General Library in the University of Tokyo"
$ws.Range("M3").Value = "394c6988-9bd0-4adc-8c3c-4e05cb02b6a2"
$ws.Range("O3").Value = "キリシタン写本"

Add-HyperlinkCell "I3" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian/document/394c6988-9bd0-4adc-8c3c-4e05cb02b6a2" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian/document/394c6988-9bd0-4adc-8c3c-4e05cb02b6a2"
Add-HyperlinkCell "J3" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/39288/manifest" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/39288/manifest"
Add-HyperlinkCell "K3" "http://iiif.io/api/presentation/2" "rightToLeftDirection" "http://iiif.io/api/presentation/2#rightToLeftDirection"
Add-HyperlinkCell "P3" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/39288" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/39288"
Add-HyperlinkCell "Q3" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/files/square/a14eada5362f9ec224a246534dd899aebf594a26.jpg" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/files/square/a14eada5362f9ec224a246534dd899aebf594a26.jpg"
Add-HyperlinkCell "R3" "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian" $null "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/christian"
Add-HyperlinkCell "S3" "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse" $null "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse"
